$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the environment-specific values from test18 -> test21
$ws.Range("A2").Value = "https://test21.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test21.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test21.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest21"
$ws.Range("G2").Value = "test21"
$ws.Range("K2").Value = "test21"

# Update the view: scroll back to A1 (remove topLeftCell="E1") and change selection to C8
[void]$ws.Range("A1").Select()
[void]$ws.Range("C8").Select()
